$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.911.84"
$ws.Range("E2").Value = "  -1.81%  "

$ws.Range("D3").Value = "1.888.63"
$ws.Range("E3").Value = "  -2.70%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.84%  "

$ws.Range("D5").Value = "'0.7341"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "

$ws.Range("D6").Value = "'242.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.39%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.62%  "

$ws.Range("E8").Value = "  -2.38%  "

$ws.Range("D9").Value = "'26.20"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.00%  "

$ws.Range("D10").Value = "'0.06901"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.22%  "

$ws.Range("D11").Value = "'0.7713"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.13%  "

$ws.Range("D12").Value = "'0.07942"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.97%  "

$ws.Range("D13").Value = "1.880.19"
$ws.Range("E13").Value = "  -3.07%  "

$ws.Range("D14").Value = "'5.223"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.40%  "

$ws.Range("D15").Value = "'91.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.29%  "

$ws.Range("D16").Value = "'14.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.66%  "

$ws.Range("D17").Value = "29.904.42"
$ws.Range("E17").Value = "  -1.78%  "

$ws.Range("D18").Value = "'5.745"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("D19").Value = "'239.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.88%  "

$ws.Range("D20").Value = "'0.000007751"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.14%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.45%  "

$ws.Range("D22").Value = "2.122.20"
$ws.Range("E22").Value = "  -3.09%  "

$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("D24").Value = "'6.905"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.54%  "

$ws.Range("D25").Value = "'9.299"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.16%  "

$ws.Range("D26").Value = "'164.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.69%  "

$ws.Range("D27").Value = "'18.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("D28").Value = "'0.1266"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.22%  "

$ws.Range("D29").Value = "'2.012"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.87%  "

$ws.Range("D30").Value = "'1.354"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.91%  "

$ws.Range("D31").Value = "'1.533"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.16%  "

$ws.Range("D32").Value = "'4.302"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.88%  "

$ws.Range("D33").Value = "'4.074"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.99%  "

$ws.Range("D34").Value = "'0.05103"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("D35").Value = "'1.279"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").Value = "'0.7357"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.43%  "

$ws.Range("D37").Value = "'2.720"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.52%  "

$ws.Range("D38").Value = "'0.01918"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.76%  "

$ws.Range("D39").Value = "'2.773"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.58%  "

$ws.Range("D40").Value = "'6.293"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.19%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.4461"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.37%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'74.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.49%  "

$ws.Range("D43").Value = "'1.931"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.63%  "

$ws.Range("E44").Value = "  -0.55%  "

$ws.Range("D45").Value = "'0.8374"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").Value = "'7.649"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.23%  "

$ws.Range("D47").Value = "'100.69"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").Value = "'9.777"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("D49").Value = "'36.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.82%  "

$ws.Range("D50").Value = "2.024.62"
$ws.Range("E50").Value = "  -2.87%  "

$ws.Range("D51").Value = "'942.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.37%  "
